# TillDateRef.xlsx - "Add files via upload" edit
#
# Rows 20 and 21 already contained the Date/Content/Reference columns
# (B/C/D) for entries #12 and #13, but the leading "S.No." column (A)
# was left blank. This fills in the missing sequence numbers, matching
# the pattern used by rows 15-19 (A15=7 ... A19=11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = 12
$ws.Range("A21").Value = 13

# Leave the sheet's active selection on the next empty row, as in the
# saved workbook.
$ws.Range("A22").Select() | Out-Null
